# EMSC Stats: add tie-breaking rule
#  - new columns Q (Final Points From) / R (Semi Points From)
#  - fix "Marocco" -> "Morocco" typo
#  - re-sort the non-finalist block (rows 27-35) by semi placement,
#    which re-ties "Moldova, Republic of" down to "Moldova"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new headers (columns Q/R) ------------------------------------------
$ws.Cells.Item(1, 17).Value = 'Final Points From'
$ws.Cells.Item(1, 18).Value = 'Semi Points From'

# match the authored column widths for the two new columns
$ws.Columns.Item(17).ColumnWidth = 17.16666666666667
$ws.Columns.Item(18).ColumnWidth = 16.16666666666667

# --- row 25: country spelling fix ---------------------------------------
$ws.Cells.Item(25, 3).Value = 'Morocco'

# --- rows 27-35: re-sorted by semi placement (tie-break rule) -----------
$ws.Cells.Item(27, 3).Value = 'Poland'
$ws.Cells.Item(27, 4).Value = 'Viki Gabor'
$ws.Cells.Item(27, 5).Value = 'Barbie'
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 10).Value = 46
$ws.Cells.Item(27, 11).Value = 'Luke'
$ws.Cells.Item(27, 14).Value = 3

$ws.Cells.Item(28, 3).Value = 'Montenegro'
$ws.Cells.Item(28, 4).Value = 'Emel'
$ws.Cells.Item(28, 5).Value = 'Gdje je'
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 10).Value = 45
$ws.Cells.Item(28, 11).Value = 'Rodrigo  Erazo'
$ws.Cells.Item(28, 14).Value = 5

$ws.Cells.Item(29, 3).Value = 'Georgia'
$ws.Cells.Item(29, 4).Value = 'Katie Melua'
$ws.Cells.Item(29, 5).Value = 'A love like that'
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 9).Value = 15
$ws.Cells.Item(29, 10).Value = 37
$ws.Cells.Item(29, 11).Value = 'Richard Cox'
$ws.Cells.Item(29, 14).Value = 8

$ws.Cells.Item(30, 3).Value = 'Ukraine'
$ws.Cells.Item(30, 4).Value = 'Ruma'
$ws.Cells.Item(30, 5).Value = '\u0414\u043e\u0442\u0438\u043a'
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 9).Value = 16
$ws.Cells.Item(30, 10).Value = 31
$ws.Cells.Item(30, 11).Value = 'Lu\u00eds Coelho'
$ws.Cells.Item(30, 14).Value = 14

$ws.Cells.Item(31, 3).Value = 'Estonia'
$ws.Cells.Item(31, 4).Value = 'Liis Lemsalu'
$ws.Cells.Item(31, 5).Value = 'Kehakeel'
$ws.Cells.Item(31, 9).Value = 13
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 'Mathias'
$ws.Cells.Item(31, 14).Value = 2

$ws.Cells.Item(32, 3).Value = 'Romania'
$ws.Cells.Item(32, 4).Value = 'WRS'
$ws.Cells.Item(32, 5).Value = 'Dale'
$ws.Cells.Item(32, 6).Value = 2
$ws.Cells.Item(32, 9).Value = 14
$ws.Cells.Item(32, 10).Value = 39
$ws.Cells.Item(32, 11).Value = 'Edu Padr\u00f3s Creus'

$ws.Cells.Item(33, 3).Value = 'Moldova'
$ws.Cells.Item(33, 4).Value = 'Vanotek & Eneli'
$ws.Cells.Item(33, 5).Value = 'Back to me'
$ws.Cells.Item(33, 6).Value = 2
$ws.Cells.Item(33, 9).Value = 15
$ws.Cells.Item(33, 10).Value = 37
$ws.Cells.Item(33, 11).Value = 'FabioMassimo'
$ws.Cells.Item(33, 14).Value = 16

$ws.Cells.Item(34, 3).Value = 'Azerbaijan'
$ws.Cells.Item(34, 4).Value = 'R\u00f6ya & Nicat R\u0259himov'
$ws.Cells.Item(34, 5).Value = 'D\u0259li kimi'
$ws.Cells.Item(34, 6).Value = 2
$ws.Cells.Item(34, 9).Value = 16
$ws.Cells.Item(34, 10).Value = 36
$ws.Cells.Item(34, 11).Value = 'Nijat'
$ws.Cells.Item(34, 14).Value = 6

$ws.Cells.Item(35, 3).Value = 'Greece'
$ws.Cells.Item(35, 4).Value = 'ZAF'
$ws.Cells.Item(35, 5).Value = 'Pes'
$ws.Cells.Item(35, 6).Value = 2
$ws.Cells.Item(35, 9).Value = 17
$ws.Cells.Item(35, 10).Value = 23
$ws.Cells.Item(35, 11).Value = 'Christoforos Andrianos'
$ws.Cells.Item(35, 14).Value = 11
